$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10; existing rows 10-64 shift down to 11-65.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new market entry.
$ws.Cells.Item(10, 1).Value = 10
$ws.Cells.Item(10, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(10, 3).Value = "La Araucanía"
$ws.Cells.Item(10, 4).Value = 44802
$ws.Cells.Item(10, 5).Value = 9
$ws.Cells.Item(10, 6).Value = "Fruta"
$ws.Cells.Item(10, 7).Value = 100108
$ws.Cells.Item(10, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(10, 9).Value = 100108003
$ws.Cells.Item(10, 10).Value = "Maracuyá"
$ws.Cells.Item(10, 11).Value = "Sin especificar"
$ws.Cells.Item(10, 12).Value = "Primera"
$ws.Cells.Item(10, 13).Value = 50
$ws.Cells.Item(10, 14).Value = 36000
$ws.Cells.Item(10, 15).Value = 36000
$ws.Cells.Item(10, 16).Value = 36000
$ws.Cells.Item(10, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(10, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(10, 19).Value = 2000
$ws.Cells.Item(10, 20).Value = 18
